$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.557.89'
$ws.Range("E2").Value = '  +4.98%  '
$ws.Range("D3").Value = '2.472.68'
$ws.Range("E3").Value = '  +6.09%  '
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '568.01'
$ws.Range("E5").Value = '  +3.88%  '
$ws.Range("D6").Value = '143.62'
$ws.Range("E6").Value = '  +9.72%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +2.54%  '
$ws.Range("D9").Value = '2.473.44'
$ws.Range("E9").Value = '  +6.18%  '
$ws.Range("D10").Value = '0.106'
$ws.Range("E10").Value = '  +4.73%  '
$ws.Range("D11").Value = '5.74'
$ws.Range("E11").Value = '  +4.22%  '
$ws.Range("E12").Value = '  +1.08%  '
$ws.Range("E13").Value = '  +5.18%  '
$ws.Range("D14").Value = '26.41'
$ws.Range("E14").Value = '  +12.23%  '
$ws.Range("D15").Value = '2.917.11'
$ws.Range("E15").Value = '  +5.90%  '
$ws.Range("D16").Value = '63.424.02'
$ws.Range("E16").Value = '  +4.70%  '
$ws.Range("D17").Value = '0.0000143'
$ws.Range("E17").Value = '  +7.01%  '
$ws.Range("D18").Value = '2.475.14'
$ws.Range("E18").Value = '  +5.37%  '
$ws.Range("D19").Value = '11.28'
$ws.Range("E19").Value = '  +6.29%  '
$ws.Range("D20").Value = '341.94'
$ws.Range("E20").Value = '  +8.77%  '
$ws.Range("D21").Value = '4.31'
$ws.Range("E21").Value = '  +5.85%  '
$ws.Range("D22").Value = '6.83'
$ws.Range("E22").Value = '  +3.95%  '
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = '65.73'
$ws.Range("E24").Value = '  +2.85%  '
$ws.Range("D25").Value = '0.175'
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("D27").Value = '1.52'
$ws.Range("E27").Value = '  +9.73%  '
$ws.Range("D28").Value = '8.19'
$ws.Range("E28").Value = '  +3.70%  '
$ws.Range("D29").Value = '1.33'
$ws.Range("E29").Value = '  +6.82%  '
$ws.Range("D30").Value = '0.0₃0824'
$ws.Range("E30").Value = '  +13.04%  '
$ws.Range("D31").Value = '6.84'
$ws.Range("E31").Value = '  +14.98%  '
$ws.Range("D32").Value = '1.86'
$ws.Range("E32").Value = '  +7.14%  '
$ws.Range("D33").Value = '177.17'
$ws.Range("E33").Value = '  +2.27%  '
$ws.Range("E34").Value = '  +10.94%  '
$ws.Range("D35").Value = '0.401'
$ws.Range("E35").Value = '  +4.41%  '
$ws.Range("D36").Value = '18.98'
$ws.Range("E36").Value = '  +5.58%  '
$ws.Range("D37").Value = '372.91'
$ws.Range("E37").Value = '  +15.78%  '
$ws.Range("D38").Value = '4.47'
$ws.Range("E38").Value = '  +8.43%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  +12.53%  '
$ws.Range("E42").Value = '  +6.42%  '
$ws.Range("D43").Value = '150.45'
$ws.Range("E43").Value = '  +9.86%  '
$ws.Range("D44").Value = '3.72'
$ws.Range("E44").Value = '  +6.49%  '
$ws.Range("D45").Value = '20.84'
$ws.Range("E45").Value = '  +9.09%  '
$ws.Range("E46").Value = '  +5.95%  '
$ws.Range("D47").Value = '0.0966'
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("D48").Value = '0.0524'
$ws.Range("E48").Value = '  +5.62%  '
$ws.Range("D49").Value = '0.0₆0236'
$ws.Range("E49").Value = '  +7.04%  '
$ws.Range("D50").Value = '0.0226'
$ws.Range("E50").Value = '  +4.78%  '
$ws.Range("D51").Value = '18.17'
$ws.Range("E51").Value = '  +7.47%  '
